$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextCell $ws 'D2' '27.577.93'
Set-TextCell $ws 'E2' '  -0.68%  '
Set-TextCell $ws 'D3' '1.748.07'
Set-TextCell $ws 'E3' '  -0.13%  '
Set-TextCell $ws 'D4' '1.002'
Set-TextCell $ws 'E4' '  +0.06%  '
Set-TextCell $ws 'D5' '323.46'
Set-TextCell $ws 'E5' '  +0.83%  '
Set-TextCell $ws 'D6' '1.001'
Set-TextCell $ws 'E6' '  +0.11%  '
Set-TextCell $ws 'D7' '0.4622'
Set-TextCell $ws 'E7' '  +9.62%  '
Set-TextCell $ws 'D8' '0.3543'
Set-TextCell $ws 'E8' '  -2.16%  '
Set-TextCell $ws 'D9' '0.07453'
Set-TextCell $ws 'E9' '  +0.32%  '
Set-TextCell $ws 'E10' '  -0.81%  '
Set-TextCell $ws 'E11' '  +0.70%  '
Set-TextCell $ws 'E12' '  +0.13%  '
Set-TextCell $ws 'E13' '  -0.07%  '
Set-TextCell $ws 'D14' '5.976'
Set-TextCell $ws 'E14' '  -0.94%  '
Set-TextCell $ws 'D15' '7.088'
Set-TextCell $ws 'E15' '  -2.39%  '
Set-TextCell $ws 'D16' '1.742.09'
Set-TextCell $ws 'E16' '  -1.42%  '
Set-TextCell $ws 'D17' '91.90'
Set-TextCell $ws 'E17' '  +1.27%  '
Set-TextCell $ws 'D18' '0.00001060'
Set-TextCell $ws 'E18' '  +0.86%  '
Set-TextCell $ws 'E19' '  +1.04%  '
Set-TextCell $ws 'D20' '1.001'
Set-TextCell $ws 'E20' '  +0.04%  '
Set-TextCell $ws 'E21' '  -1.60%  '
Set-TextCell $ws 'D22' '5.765'
Set-TextCell $ws 'E22' '  -2.36%  '
Set-TextCell $ws 'D23' '27.635.70'
Set-TextCell $ws 'E24' '  +0.28%  '
Set-TextCell $ws 'D25' '2.109'
Set-TextCell $ws 'E25' '  +0.39%  '
Set-TextCell $ws 'D26' '163.05'
Set-TextCell $ws 'E26' '  +3.85%  '
Set-TextCell $ws 'D27' '20.08'
Set-TextCell $ws 'E27' '  -0.02%  '
Set-TextCell $ws 'D28' '1.946.62'
Set-TextCell $ws 'E28' '  -0.87%  '
Set-TextCell $ws 'B29' 'BitcoinCash'
Set-TextCell $ws 'C29' 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextCell $ws 'D29' '125.78'
Set-TextCell $ws 'E29' '  +1.85%  '
Set-TextCell $ws 'B30' 'LidoDAOToken'
Set-TextCell $ws 'C30' 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextCell $ws 'D30' '2.055'
Set-TextCell $ws 'E30' '  -2.99%  '
Set-TextCell $ws 'D31' '1.051'
Set-TextCell $ws 'E31' '  -5.85%  '
Set-TextCell $ws 'D32' '0.09228'
Set-TextCell $ws 'E32' '  +4.42%  '
Set-TextCell $ws 'D33' '3.667'
Set-TextCell $ws 'E33' '  +0.74%  '
Set-TextCell $ws 'D34' '5.520'
Set-TextCell $ws 'E34' '  -0.30%  '
Set-TextCell $ws 'D35' '0.02287'
Set-TextCell $ws 'E35' '  +0.37%  '
Set-TextCell $ws 'D36' '11.73'
Set-TextCell $ws 'E36' '  -4.07%  '
Set-TextCell $ws 'D37' '0.06024'
Set-TextCell $ws 'E37' '  +0.09%  '
Set-TextCell $ws 'D38' '0.2083'
Set-TextCell $ws 'E38' '  -0.23%  '
Set-TextCell $ws 'D39' '4.955'
Set-TextCell $ws 'E39' '  +0.67%  '
Set-TextCell $ws 'D40' '0.6283'
Set-TextCell $ws 'E40' '  +0.10%  '
Set-TextCell $ws 'D41' '1.195'
Set-TextCell $ws 'E41' '  +1.94%  '
Set-TextCell $ws 'D42' '1.377'
Set-TextCell $ws 'E42' '  -1.24%  '
Set-TextCell $ws 'D43' '7.735'
Set-TextCell $ws 'E43' '  -1.25%  '
Set-TextCell $ws 'D44' '13.20'
Set-TextCell $ws 'E44' '  -0.96%  '
Set-TextCell $ws 'D45' '3.709'
Set-TextCell $ws 'E45' '  +0.96%  '
Set-TextCell $ws 'E46' '  +0.27%  '
Set-TextCell $ws 'E47' '  +0.16%  '
Set-TextCell $ws 'D48' '1.933'
Set-TextCell $ws 'E48' '  -1.60%  '
Set-TextCell $ws 'B49' 'EOS'
Set-TextCell $ws 'C49' 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
Set-TextCell $ws 'D49' '1.131'
Set-TextCell $ws 'E49' '  -3.80%  '
Set-TextCell $ws 'B50' 'Cronos'
Set-TextCell $ws 'C50' 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextCell $ws 'D50' '0.06857'
Set-TextCell $ws 'E50' '  +0.95%  '
Set-TextCell $ws 'D51' '71.60'
Set-TextCell $ws 'E51' '  -2.26%  '

Write-Host "Applied 98 cell updates"